$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b4740fa38aca8f3ddf3aea9c446fa3a03106ff2c/e2e/"

# ---------------------------------------------------------------------------
# Overview sheet (row 4: 3fd41955..., row 5: 61d53267...)
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A4").Value = "3fd41955-8d61-47bc-85b7-76c3ec97a3f5.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B4"), $ghBase + "3fd41955-8d61-47bc-85b7-76c3ec97a3f5.md", "", "", "e2e\3fd41955-8d61-47bc-85b7-76c3ec97a3f5.md")
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("D4").Value = ""
$wsOv.Range("E4").Value = "Ready for handoff"
$wsOv.Range("F4").Value = "Ready for handoff"
$wsOv.Range("G4").Value = "2016-09-06 10:01:32"
$wsOv.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Range("A5").Value = "61d53267-352e-4f3c-9bd2-f8f8055bd5e0.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B5"), $ghBase + "61d53267-352e-4f3c-9bd2-f8f8055bd5e0.md", "", "", "e2e\61d53267-352e-4f3c-9bd2-f8f8055bd5e0.md")
$wsOv.Range("C5").Value = ".md"
$wsOv.Range("D5").Value = ""
$wsOv.Range("E5").Value = "Ready for handoff"
$wsOv.Range("F5").Value = "Ready for handoff"
$wsOv.Range("G5").Value = "2016-09-06 10:01:32"
$wsOv.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Columns.Item(5).ColumnWidth = 17.22
$wsOv.Columns.Item(6).ColumnWidth = 17.22

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G5"))

# ---------------------------------------------------------------------------
# zh-cn sheet (row 4: 3fd41955..., row 5: 61d53267...)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $ghBase + "3fd41955-8d61-47bc-85b7-76c3ec97a3f5.md", "", "", "3fd41955-8d61-47bc-85b7-76c3ec97a3f5.md")
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'False"
$wsZh.Range("G4").Value = "3fd41955-8d61-47bc-85b7-76c3ec97a3f5.60fefb14d731d9fee4b554623e6f4f5b8bbde58a.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-06 10:01:11"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("O4").Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $ghBase + "61d53267-352e-4f3c-9bd2-f8f8055bd5e0.md", "", "", "61d53267-352e-4f3c-9bd2-f8f8055bd5e0.md")
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = "61d53267-352e-4f3c-9bd2-f8f8055bd5e0.5c1214d0a46f08716ce99a4070abdd0a618eccdd.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-09-06 10:01:11"
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M5").Value = "'True"
$wsZh.Range("O5").Value = "'False"

$wsZh.Columns.Item(3).ColumnWidth = 17.22

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# de-de sheet (row 4: 3fd41955..., row 5: 61d53267...)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $ghBase + "3fd41955-8d61-47bc-85b7-76c3ec97a3f5.md", "", "", "3fd41955-8d61-47bc-85b7-76c3ec97a3f5.md")
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'False"
$wsDe.Range("G4").Value = "3fd41955-8d61-47bc-85b7-76c3ec97a3f5.60fefb14d731d9fee4b554623e6f4f5b8bbde58a.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-06 10:01:32"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("O4").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $ghBase + "61d53267-352e-4f3c-9bd2-f8f8055bd5e0.md", "", "", "61d53267-352e-4f3c-9bd2-f8f8055bd5e0.md")
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = "61d53267-352e-4f3c-9bd2-f8f8055bd5e0.5c1214d0a46f08716ce99a4070abdd0a618eccdd.de-de.xlf"
$wsDe.Range("H5").Value = "2016-09-06 10:01:32"
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M5").Value = "'True"
$wsDe.Range("O5").Value = "'False"

$wsDe.Columns.Item(3).ColumnWidth = 17.22

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))
